# Add a new driver ("prakash") as the first row of data on the "driver"
# sheet, pushing the existing drivers (umesh, pradeep, paramesh, chandru)
# down by one row, and renumber the driver_id column 1..5. Make the
# "driver" sheet the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("driver")

# Selecting/activating the "driver" tab makes it the active sheet, which
# flips tabSelected on the sheetViews and sets workbookView's activeTab.
$ws.Activate()

# Existing driver names (read before overwriting) so we can shift them
# down by one row. Note: the reflection-based COM bridge needs the
# getter invoked explicitly (Value()) -- bare ".Value" yields the
# property descriptor, not the cell's contents.
$row2Name = $ws.Range("B2").Value()
$row3Name = $ws.Range("B3").Value()
$row4Name = $ws.Range("B4").Value()
$row5Name = $ws.Range("B5").Value()

# New first driver entry.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "prakash"
$ws.Range("C2").Value = 1234

# Shift the pre-existing drivers down one row each, renumbering driver_id.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = $row2Name
$ws.Range("C3").Value = 1234

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = $row3Name
$ws.Range("C4").Value = 1234

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = $row4Name
$ws.Range("C5").Value = 1234

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = $row5Name
$ws.Range("C6").Value = 1234

# Update the saved selection on the driver sheet.
$ws.Range("C3").Select()
